$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "Hello"
